# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values (header at G1 = "K").
# Update the literal K values for rows 2-9 per the regenerated save data.
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 3
$ws.Range("G9").Value = 1
